$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: volume/number and week-covering dates ---
$ws.Range("A8").Value = "Volume 32   Number  50"
$ws.Range("C9").Value = "Report Covering the Week  12/8/2025  Through  12/14/2025"

# --- Type-changing cells (text <-> number) ---
# Row 15: G15 (1 -> "0"), H15 (-100 -> "***.*")
$ws.Range("C15").Copy($ws.Range("G15"))
$ws.Range("E15").Copy($ws.Range("H15"))

# Row 16: C16 (1 -> "0")
$ws.Range("C15").Copy($ws.Range("C16"))

# Row 18: D18 ("0" -> 1), E18 ("***.*" -> 100)
$ws.Range("C17").Copy($ws.Range("D18"))
$ws.Range("D18").Value = 1
$ws.Range("K15").Copy($ws.Range("E18"))
$ws.Range("E18").Value = 100

# Row 20: C20 ("0" -> 10)
$ws.Range("G20").Copy($ws.Range("C20"))
$ws.Range("C20").Value = 10

# Row 25: C25 (2 -> "0")
$ws.Range("C23").Copy($ws.Range("C25"))

# Row 27: G27 (1 -> "0"), H27 (-100 -> "***.*")
$ws.Range("C27").Copy($ws.Range("G27"))
$ws.Range("E27").Copy($ws.Range("H27"))

# Row 28: D28 ("0" -> 1), E28 ("***.*" -> -100)
$ws.Range("I28").Copy($ws.Range("D28"))
$ws.Range("D28").Value = 1
$ws.Range("K27").Copy($ws.Range("E28"))
$ws.Range("E28").Value = -100

# Row 29: D29 (2 -> "0"), E29 (-100 -> "***.*")
$ws.Range("C29").Copy($ws.Range("D29"))
$ws.Range("M29").Copy($ws.Range("E29"))

# Row 30: D30 (2 -> "0"), E30 (-100 -> "***.*")
$ws.Range("C30").Copy($ws.Range("D30"))
$ws.Range("M30").Copy($ws.Range("E30"))

# --- Plain numeric value updates ---
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 4
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = -20
$ws.Range("J16").Value = 77
$ws.Range("K16").Value = -2.597402597402
$ws.Range("L16").Value = -7.407407407407

$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 33.333333333333
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = -33.333333333333
$ws.Range("I17").Value = 246
$ws.Range("J17").Value = 220
$ws.Range("K17").Value = 11.818181818181
$ws.Range("L17").Value = 14.418604651162

$ws.Range("C18").Value = 2
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = -40
$ws.Range("I18").Value = 57
$ws.Range("J18").Value = 68
$ws.Range("K18").Value = -16.176470588235
$ws.Range("L18").Value = -20.833333333333

$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 24
$ws.Range("G19").Value = 12
$ws.Range("H19").Value = 100
$ws.Range("I19").Value = 237
$ws.Range("J19").Value = 259
$ws.Range("K19").Value = -8.494208494208
$ws.Range("L19").Value = -19.112627986348

$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 17
$ws.Range("H20").Value = -10.526315789473
$ws.Range("I20").Value = 165
$ws.Range("J20").Value = 182
$ws.Range("K20").Value = -9.340659340659
$ws.Range("L20").Value = -19.512195121951

$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 14
$ws.Range("E21").Value = 57.142857142857
$ws.Range("F21").Value = 62
$ws.Range("G21").Value = 62
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 807
$ws.Range("J21").Value = 823
$ws.Range("K21").Value = -1.944106925880
$ws.Range("L21").Value = -8.399545970488

$ws.Range("D24").Value = 17
$ws.Range("E24").Value = -29.411764705882
$ws.Range("G24").Value = 50
$ws.Range("H24").Value = 22
$ws.Range("I24").Value = 574
$ws.Range("J24").Value = 561
$ws.Range("K24").Value = 2.317290552584
$ws.Range("L24").Value = -1.204819277108

$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -100
$ws.Range("F25").Value = 17
$ws.Range("G25").Value = 18
$ws.Range("H25").Value = -5.555555555555
$ws.Range("J25").Value = 157
$ws.Range("K25").Value = -29.299363057324
$ws.Range("L25").Value = -29.299363057324

$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 27
$ws.Range("G26").Value = 46
$ws.Range("H26").Value = -41.304347826087
$ws.Range("I26").Value = 436
$ws.Range("J26").Value = 434
$ws.Range("K26").Value = 0.460829493087
$ws.Range("L26").Value = 14.736842105263

$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 24
$ws.Range("K28").Value = -8.333333333333
$ws.Range("L28").Value = -15.384615384615

$ws.Range("I33").Value = 2
